# repull data, push all data, mean calculation
# Update column F (dSF) values to reflect the re-pulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = -1
$ws.Range("F5").Value = -11
$ws.Range("F6").Value = -5
$ws.Range("F7").Value = -1
$ws.Range("F8").Value = -13
$ws.Range("F9").Value = 8
$ws.Range("F11").Value = -3
$ws.Range("F12").Value = -1
$ws.Range("F13").Value = -2
